$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Insert a new row at position 16 (shifts old rows 16-23 down to 17-24)
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new job listing
$ws.Range("A16").Value = "2026-02-03 18:54:00"
$ws.Range("B16").Value = 'ファイルメーカーでの在庫・顧客管理システム構築'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5485054')
$ws.Range("F16").Style = "Hyperlink"
$ws.Range("G16").Value = 53
$ws.Range("H16").Value = '◇管理'

# Refresh the acquisition timestamp (column A) for every other listing row
for ($r = 2; $r -le 24; $r++) {
    if ($r -ne 16) {
        $ws.Range("A" + $r).Value = "2026-02-03 18:54:00"
    }
}
